$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (single-decimal-point values),
# so they stay text exactly like the other price cells in column D.
$textCells = @('D5', 'D6', 'D9', 'D10', 'D14', 'D15', 'D18', 'D21', 'D23', 'D24', 'D26', 'D28', 'D32', 'D33', 'D34', 'D36', 'D38', 'D39', 'D42', 'D43', 'D44', 'D45', 'D48', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values from the crypto-price refresh.
$ws.Range('D2').Value = '51.557.34'
$ws.Range('E2').Value = '  +1.50%  '
$ws.Range('D3').Value = '2.990.56'
$ws.Range('E3').Value = '  +2.24%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '381.97'
$ws.Range('E5').Value = '  +1.85%  '
$ws.Range('D6').Value = '104.27'
$ws.Range('E6').Value = '  +4.32%  '
$ws.Range('E7').Value = '  +1.90%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = '0.597'
$ws.Range('E9').Value = '  +2.29%  '
$ws.Range('D10').Value = '36.83'
$ws.Range('E10').Value = '  +2.51%  '
$ws.Range('E12').Value = '  +2.09%  '
$ws.Range('D13').Value = '3.467.87'
$ws.Range('E13').Value = '  +2.34%  '
$ws.Range('D14').Value = '18.53'
$ws.Range('E14').Value = '  +3.19%  '
$ws.Range('D15').Value = '7.84'
$ws.Range('E15').Value = '  +3.63%  '
$ws.Range('D16').Value = '2.980.06'
$ws.Range('E16').Value = '  +2.09%  '
$ws.Range('E17').Value = '  -0.31%  '
$ws.Range('D18').Value = '0.996'
$ws.Range('E18').Value = '  +0.79%  '
$ws.Range('D19').Value = '51.643.23'
$ws.Range('E19').Value = '  +1.75%  '
$ws.Range('E20').Value = '  +0.85%  '
$ws.Range('D21').Value = '12.56'
$ws.Range('E21').Value = '  +1.58%  '
$ws.Range('D22').Value = '0.0₃0965'
$ws.Range('D23').Value = '70.42'
$ws.Range('E23').Value = '  +2.35%  '
$ws.Range('D24').Value = '267.63'
$ws.Range('E24').Value = '  +1.08%  '
$ws.Range('E25').Value = '  +2.20%  '
$ws.Range('D26').Value = '8.07'
$ws.Range('E26').Value = '  +0.89%  '
$ws.Range('E27').Value = '  +4.51%  '
$ws.Range('D28').Value = '7.21'
$ws.Range('E28').Value = '  -2.63%  '
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('E30').Value = '  +2.67%  '
$ws.Range('E31').Value = '  +0.26%  '
$ws.Range('D32').Value = '10.40'
$ws.Range('E32').Value = '  +4.38%  '
$ws.Range('D33').Value = '34.72'
$ws.Range('E33').Value = '  +5.03%  '
$ws.Range('D34').Value = '51.41'
$ws.Range('E34').Value = '  +1.12%  '
$ws.Range('E35').Value = '  +0.62%  '
$ws.Range('D36').Value = '0.0445'
$ws.Range('E36').Value = '  +2.02%  '
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('D38').Value = '3.29'
$ws.Range('E38').Value = '  +6.84%  '
$ws.Range('D39').Value = '17.04'
$ws.Range('E39').Value = '  +3.78%  '
$ws.Range('E40').Value = '  +5.33%  '
$ws.Range('D42').Value = '1.85'
$ws.Range('E42').Value = '  +2.01%  '
$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D43').Value = '127.48'
$ws.Range('E43').Value = '  +6.81%  '
$ws.Range('B44').Value = 'NEARProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D44').Value = '3.86'
$ws.Range('E44').Value = '  +14.39%  '
$ws.Range('D45').Value = '21.43'
$ws.Range('E45').Value = '  +2.10%  '
$ws.Range('E46').Value = '  +0.43%  '
$ws.Range('E47').Value = '  +2.26%  '
$ws.Range('D48').Value = '2.36'
$ws.Range('E48').Value = '  +1.03%  '
$ws.Range('D49').Value = '2.037.83'
$ws.Range('E49').Value = '  +2.59%  '
$ws.Range('D50').Value = '3.285.42'
$ws.Range('E50').Value = '  +2.28%  '
$ws.Range('B51').Value = 'BEAM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range('D51').Value = '0.0331'
$ws.Range('E51').Value = '  +2.57%  '
